$wb = $excel.ActiveWorkbook

# --- Order Details sheet: fill in column F (Test Result) with TRUE for rows 2-41 ---
$wsOrderDetails = $wb.Worksheets.Item("Order Details")
$wsOrderDetails.Range("F2:F41").Value = $true

# Make Order Details the active sheet, scrolled to show row 32 at top
$wsOrderDetails.Activate()
$excel.ActiveWindow.ScrollRow = 32
$wsOrderDetails.Range("F1").Select()

# --- Products sheet: no longer the active/selected tab (keep its scroll position) ---
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("E1").Select()
